$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 220
$ws.Range("J9").Value = 220
$ws.Range("L9").Value = 220
$ws.Range("N9").Value = -558
$ws.Range("H125").Value = 3636
$ws.Range("J125").Value = 3636
$ws.Range("L125").Value = 32724
$ws.Range("N125").Value = -37644
$ws.Range("H132").Value = 1198.8043
$ws.Range("I132").Value = 1116.9546
$ws.Range("K132").Value = 3350.8638
$ws.Range("M132").Value = -820.8638000000001
$ws.Range("H137").Value = 2375.077
$ws.Range("I137").Value = 1798.6666
$ws.Range("K137").Value = 5395.9998
$ws.Range("M137").Value = -2845.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1491.7059
$ws.Range("I45").Value = 911.6667
$ws.Range("K45").Value = 911.6667
$ws.Range("M45").Value = -534.6667
$ws.Range("H61").Value = 4948.8335
$ws.Range("I61").Value = 2223.5
$ws.Range("K61").Value = 2223.5
$ws.Range("M61").Value = -2011.5
$ws.Range("H63").Value = 5849.6665
$ws.Range("I63").Value = 6699.8
$ws.Range("K63").Value = 6699.8
$ws.Range("M63").Value = -6013.8
$ws.Range("H66").Value = 5849.6665
$ws.Range("I66").Value = 6699.8
$ws.Range("K66").Value = 33499
$ws.Range("M66").Value = -30067
$ws.Range("H74").Value = 3289.6365
$ws.Range("I74").Value = 3100
$ws.Range("K74").Value = 3100
$ws.Range("M74").Value = -2226
$ws.Range("H77").Value = 3289.6365
$ws.Range("I77").Value = 3100
$ws.Range("K77").Value = 15500
$ws.Range("M77").Value = -11132
$ws.Range("H88").Value = 3032.818
$ws.Range("I88").Value = 1538.75
$ws.Range("J88").Value = 3886.5715
$ws.Range("K88").Value = 1538.75
$ws.Range("L88").Value = 3886.5715
$ws.Range("M88").Value = -1132.75
$ws.Range("N88").Value = -4698.5715
$ws.Range("H91").Value = 3032.818
$ws.Range("I91").Value = 1538.75
$ws.Range("J91").Value = 3886.5715
$ws.Range("K91").Value = 1538.75
$ws.Range("L91").Value = 3886.5715
$ws.Range("M91").Value = -134.75
$ws.Range("N91").Value = -6694.5715
$ws.Range("H110").Value = 244.8
$ws.Range("I110").Value = 224.375
$ws.Range("K110").Value = 224.375
$ws.Range("M110").Value = 1820.625
$ws.Range("H122").Value = 1579.3077
$ws.Range("I122").Value = 1536.8572
$ws.Range("K122").Value = 4610.571599999999
$ws.Range("M122").Value = -2160.571599999999
$ws.Range("H136").Value = 4948.8335
$ws.Range("I136").Value = 2223.5
$ws.Range("K136").Value = 6670.5
$ws.Range("M136").Value = -4120.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 631.4545000000001
$ws.Range("I107").Value = 468.25
$ws.Range("K107").Value = 468.25
$ws.Range("M107").Value = 1451.75
$ws.Range("H134").Value = 6483.826
$ws.Range("I134").Value = 6864.4287
$ws.Range("K134").Value = 20593.2861
$ws.Range("M134").Value = -18058.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2564.35
$ws.Range("I31").Value = 2326.3333
$ws.Range("J31").Value = 2921.375
$ws.Range("K31").Value = 2326.3333
$ws.Range("L31").Value = 2921.375
$ws.Range("M31").Value = -2031.3333
$ws.Range("N31").Value = -3511.375
$ws.Range("H34").Value = 2564.35
$ws.Range("I34").Value = 2326.3333
$ws.Range("J34").Value = 2921.375
$ws.Range("K34").Value = 2326.3333
$ws.Range("L34").Value = 2921.375
$ws.Range("M34").Value = -2124.3333
$ws.Range("N34").Value = -3325.375
$ws.Range("H58").Value = 3892.0625
$ws.Range("I58").Value = 2959.4
$ws.Range("J58").Value = 4316
$ws.Range("K58").Value = 2959.4
$ws.Range("L58").Value = 4316
$ws.Range("M58").Value = -2756.4
$ws.Range("N58").Value = -4722
$ws.Range("H122").Value = 4535.2856
$ws.Range("I122").Value = 3373
$ws.Range("J122").Value = 6085
$ws.Range("K122").Value = 10119
$ws.Range("L122").Value = 18255
$ws.Range("M122").Value = -7669
$ws.Range("N122").Value = -23155
$ws.Range("H134").Value = 1873.9231
$ws.Range("I134").Value = 1705.7273
$ws.Range("K134").Value = 5117.1819
$ws.Range("M134").Value = -2582.1819
$ws.Range("H136").Value = 3892.0625
$ws.Range("I136").Value = 2959.4
$ws.Range("J136").Value = 4316
$ws.Range("K136").Value = 8878.200000000001
$ws.Range("L136").Value = 12948
$ws.Range("M136").Value = -6328.200000000001
$ws.Range("N136").Value = -18048
$ws.Range("H138").Value = 105803.336
$ws.Range("J138").Value = 105803.336
$ws.Range("L138").Value = 105803.336
$ws.Range("N138").Value = -116083.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2050.2727
$ws.Range("J122").Value = 2050.2727
$ws.Range("L122").Value = 18452.4543
$ws.Range("N122").Value = -23352.4543
$ws.Range("H132").Value = 1421.6666
$ws.Range("J132").Value = 1799.091
$ws.Range("L132").Value = 16191.819
$ws.Range("N132").Value = -21251.819
$ws.Range("H139").Value = 8193.4375
$ws.Range("I139").Value = 8193.4375
$ws.Range("K139").Value = 24580.3125
$ws.Range("M139").Value = -19440.3125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 3675157.8
$ws.Range("J7").Value = 761142.5600000001
$ws.Range("L7").Value = 761142.5600000001
$ws.Range("N7").Value = -761366.5600000001
$ws.Range("H8").Value = 3675157.8
$ws.Range("J8").Value = 761142.5600000001
$ws.Range("L8").Value = 761142.5600000001
$ws.Range("N8").Value = -761420.5600000001
$ws.Range("H97").Value = 2327.6
$ws.Range("I97").Value = 2182.7144
$ws.Range("J97").Value = 2665.6667
$ws.Range("K97").Value = 2182.7144
$ws.Range("L97").Value = 2665.6667
$ws.Range("M97").Value = -1686.7144
$ws.Range("N97").Value = -3657.6667
$ws.Range("H102").Value = 3677.1035
$ws.Range("I102").Value = 4264.6113
$ws.Range("K102").Value = 4264.6113
$ws.Range("M102").Value = -2642.6113
$ws.Range("H122").Value = 1365.6364
$ws.Range("I122").Value = 1365.6364
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4096.9092
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1646.9092
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 26933.83
$ws.Range("I126").Value = 2790.6
$ws.Range("J126").Value = 40862.617
$ws.Range("K126").Value = 8371.799999999999
$ws.Range("L126").Value = 122587.851
$ws.Range("M126").Value = -5901.799999999999
$ws.Range("N126").Value = -127527.851
$ws.Range("H132").Value = 1889.0541
$ws.Range("I132").Value = 1427.7916
$ws.Range("K132").Value = 4283.3748
$ws.Range("M132").Value = -1753.3748

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1087.1428
$ws.Range("I22").Value = 602.25
$ws.Range("J22").Value = 1201.2354
$ws.Range("K22").Value = 602.25
$ws.Range("L22").Value = 1201.2354
$ws.Range("M22").Value = -307.25
$ws.Range("N22").Value = -1791.2354
$ws.Range("H27").Value = 1087.1428
$ws.Range("I27").Value = 602.25
$ws.Range("J27").Value = 1201.2354
$ws.Range("K27").Value = 602.25
$ws.Range("L27").Value = 1201.2354
$ws.Range("M27").Value = -495.25
$ws.Range("N27").Value = -1415.2354
$ws.Range("H61").Value = 1759.6875
$ws.Range("I61").Value = 1222.1666
$ws.Range("K61").Value = 1222.1666
$ws.Range("M61").Value = -1020.1666
$ws.Range("H82").Value = 1851.3334
$ws.Range("I82").Value = 1241.2
$ws.Range("J82").Value = 3071.6
$ws.Range("K82").Value = 1241.2
$ws.Range("L82").Value = 3071.6
$ws.Range("M82").Value = -880.2
$ws.Range("N82").Value = -3793.6
$ws.Range("H85").Value = 1851.3334
$ws.Range("I85").Value = 1241.2
$ws.Range("J85").Value = 3071.6
$ws.Range("K85").Value = 1241.2
$ws.Range("L85").Value = 3071.6
$ws.Range("M85").Value = 6.799999999999955
$ws.Range("N85").Value = -5567.6
$ws.Range("H93").Value = 15152263
$ws.Range("I93").Value = 743.8125
$ws.Range("J93").Value = 55556316
$ws.Range("K93").Value = 743.8125
$ws.Range("L93").Value = 55556316
$ws.Range("M93").Value = 504.1875
$ws.Range("N93").Value = -55558812
$ws.Range("H113").Value = 1759.6875
$ws.Range("I113").Value = 1222.1666
$ws.Range("K113").Value = 1222.1666
$ws.Range("M113").Value = 947.8334
$ws.Range("H122").Value = 4388.077
$ws.Range("I122").Value = 3920.4167
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 11761.2501
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -9311.250100000001
$ws.Range("N122").Value = -34900
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1245.4286
$ws.Range("I81").Value = 610.4
$ws.Range("K81").Value = 1220.8
$ws.Range("M81").Value = -159.8
$ws.Range("H84").Value = 1245.4286
$ws.Range("I84").Value = 610.4
$ws.Range("K84").Value = 6104
$ws.Range("M84").Value = -800
$ws.Range("H122").Value = 79296.5
$ws.Range("I122").Value = 130877.664
$ws.Range("K122").Value = 392632.992
$ws.Range("M122").Value = -390182.992
